$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data scraped on Sat Sep 23 20:36:22 UTC 2023 (GitHub Actions run)
# All Price/Volume/Coin/Link columns are plain text cells in this sheet; force text format
# on column D (Price) edits so strings like "91.00" or "7.20" keep their trailing zeros
# instead of being auto-coerced to numbers by Excel.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.700.96'
$ws.Range('E2').Value = '  +0.41%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.599.75'
$ws.Range('E3').Value = '  +0.32%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.54'
$ws.Range('E5').Value = '  +0.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.512'
$ws.Range('E6').Value = '  -0.39%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('E8').Value = '  +0.51%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.247'
$ws.Range('E9').Value = '  +1.29%  '
$ws.Range('E10').Value = '  +0.64%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0842'
$ws.Range('E11').Value = '  +0.96%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.824.24'
$ws.Range('E12').Value = '  +0.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.604.27'
$ws.Range('E13').Value = '  +0.91%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.05'
$ws.Range('E14').Value = '  +0.79%  '
$ws.Range('E15').Value = '  +0.78%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.27'
$ws.Range('E16').Value = '  +1.43%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.680.98'
$ws.Range('E17').Value = '  +0.34%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0755'
$ws.Range('E18').Value = '  +3.61%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '209.87'
$ws.Range('E19').Value = '  +1.06%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.20'
$ws.Range('E20').Value = '  +4.53%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.00'
$ws.Range('E21').Value = '  +0.16%  '
$ws.Range('E22').Value = '  +0.92%  '
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('E24').Value = '  +1.08%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.04'
$ws.Range('E25').Value = '  -1.54%  '
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.114'
$ws.Range('E28').Value = '  +0.36%  '
$ws.Range('E29').Value = '  +0.92%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0516'
$ws.Range('E30').Value = '  +2.67%  '
$ws.Range('E31').Value = '  +0.10%  '
$ws.Range('E32').Value = '  +1.10%  '
$ws.Range('E33').Value = '  +1.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.291.10'
$ws.Range('E34').Value = '  +0.86%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.620'
$ws.Range('E35').Value = '  -4.69%  '
$ws.Range('E36').Value = '  +0.81%  '
$ws.Range('E37').Value = '  +0.58%  '
$ws.Range('E38').Value = '  +0.21%  '
$ws.Range('E39').Value = '  +15.65%  '
$ws.Range('E40').Value = '  -1.96%  '
$ws.Range('E41').Value = '  -0.74%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.785'
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.18'
$ws.Range('E43').Value = '  -0.64%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.18'
$ws.Range('E44').Value = '  -0.82%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.737.86'
$ws.Range('E45').Value = '  +0.36%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.00'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.58'
$ws.Range('E47').Value = '  -0.50%  '
$ws.Range('E48').Value = '  -2.01%  '
$ws.Range('E50').Value = '  +0.27%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.39'
$ws.Range('E51').Value = '  -1.03%  '
